# Rename two table placeholder shapes on slide 2.
#   Table 10 (shape id 10) -> PH_takers_right
#   Table 5  (shape id 59) -> PH_takers_left

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)

    if ($sh.Id -eq 10 -or $sh.Name -eq "Table 10") {
        $sh.Name = "PH_takers_right"
    }
    elseif ($sh.Id -eq 59 -or $sh.Name -eq "Table 5") {
        $sh.Name = "PH_takers_left"
    }
}
